$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.172.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.241.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.576.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.850"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.246.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.019.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0982"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.77%  "
$ws.Range("E27").Value = "  -5.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.18%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0824"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("E34").Value = "  -5.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0300"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.54%  "
$ws.Range("E43").Value = "  -8.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.49%  "
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.431"
$ws.Range("D51").Style = "Normal"
